$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.712.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.617.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.72%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.605'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.82%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +13.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.63%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.62%  '

$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +14.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.005.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.25%  '

$ws.Range("E15").Value = '  +1.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.607.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.930'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.69%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '46.783.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.91%  '

$ws.Range("E20").Value = '  +7.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +31.96%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.46%  '

$ws.Range("E31").Value = '  +3.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '39.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.11%  '

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.19%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0848'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.89%  '

$ws.Range("E39").Value = '  +7.20%  '

$ws.Range("E40").Value = '  +6.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +40.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0333'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.122.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '93.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.76%  '

$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.41%  '
